$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores numeric-looking values as text (some are
# European-style grouped numbers like "66.295.11" which are not valid
# numbers at all, others are plain decimals like "8.00" that Excel would
# otherwise auto-convert to a number and silently drop the trailing
# zero / switch to scientific notation). Force the whole column to stay
# text before writing the refreshed values so every write keeps its
# literal textual representation.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "66.313.30"
$ws.Range("E2").Value = "  -1.19%  "
$ws.Range("D3").Value = "3.525.08"
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "606.85"
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("D6").Value = "145.31"
$ws.Range("E6").Value = "  -2.12%  "
$ws.Range("D7").Value = "3.524.01"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  -0.61%  "
$ws.Range("E10").Value = "  -4.39%  "
$ws.Range("D11").Value = "8.00"
$ws.Range("E11").Value = "  +1.98%  "
$ws.Range("D12").Value = "0.414"
$ws.Range("E12").Value = "  -2.43%  "
$ws.Range("D13").Value = "4.116.97"
$ws.Range("E13").Value = "  +0.07%  "
$ws.Range("D14").Value = "0.0000208"
$ws.Range("E14").Value = "  -3.77%  "
$ws.Range("D15").Value = "30.45"
$ws.Range("E15").Value = "  -4.03%  "
$ws.Range("D16").Value = "3.518.55"
$ws.Range("E16").Value = "  +0.19%  "
$ws.Range("D17").Value = "66.347.06"
$ws.Range("E17").Value = "  -1.31%  "
$ws.Range("E18").Value = "  -0.09%  "
$ws.Range("D19").Value = "10.69"
$ws.Range("E19").Value = "  -0.63%  "
$ws.Range("D20").Value = "6.21"
$ws.Range("E20").Value = "  -3.15%  "
$ws.Range("D21").Value = "14.95"
$ws.Range("E21").Value = "  -3.07%  "
$ws.Range("D22").Value = "427.29"
$ws.Range("E22").Value = "  -2.12%  "
$ws.Range("D23").Value = "0.601"
$ws.Range("E23").Value = "  -1.93%  "
$ws.Range("D24").Value = "78.22"
$ws.Range("E24").Value = "  -2.35%  "
$ws.Range("D25").Value = "3.656.28"
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("E26").Value = "  -0.49%  "
$ws.Range("D27").Value = "0.0000120"
$ws.Range("E27").Value = "  -1.16%  "
$ws.Range("D28").Value = "9.30"
$ws.Range("E28").Value = "  -5.70%  "
$ws.Range("D29").Value = "8.05"
$ws.Range("E29").Value = "  -3.56%  "
$ws.Range("E30").Value = "  -1.39%  "
$ws.Range("E31").Value = "  -0.17%  "
$ws.Range("E32").Value = "  -1.59%  "
$ws.Range("E33").Value = "  -7.55%  "
$ws.Range("D34").Value = "25.29"
$ws.Range("E34").Value = "  -0.54%  "
$ws.Range("D35").Value = "3.503.14"
$ws.Range("E35").Value = "  -0.24%  "
$ws.Range("E37").Value = "  -3.97%  "
$ws.Range("D38").Value = "7.80"
$ws.Range("E38").Value = "  -2.94%  "
$ws.Range("D39").Value = "5.62"
$ws.Range("E39").Value = "  -5.20%  "
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("D41").Value = "170.75"
$ws.Range("E41").Value = "  +0.53%  "
$ws.Range("D42").Value = "0.0860"
$ws.Range("E42").Value = "  -3.97%  "
$ws.Range("D43").Value = "5.20"
$ws.Range("E43").Value = "  -4.64%  "
$ws.Range("D44").Value = "0.890"
$ws.Range("E44").Value = "  -0.87%  "
$ws.Range("E45").Value = "  -9.59%  "
$ws.Range("D46").Value = "45.48"
$ws.Range("E46").Value = "  -0.67%  "
$ws.Range("E47").Value = "  -8.57%  "
$ws.Range("D48").Value = "25.92"
$ws.Range("E48").Value = "  -11.00%  "
$ws.Range("E49").Value = "  -1.12%  "
$ws.Range("D50").Value = "7.19"
$ws.Range("E50").Value = "  -4.00%  "
$ws.Range("D51").Value = "0.950"
$ws.Range("E51").Value = "  -3.85%  "
